# Add 2022-Q4 data: a new per-quarter sheet plus a new summary row on "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new worksheet for "2022-Q4" right after "总计" (i.e. right
#    before the current "2022-Q3" sheet), so the tab order becomes:
#    总计, 2022-Q4, 2022-Q3, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q2, 2021-Q1, 2020-Q4
# ---------------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item("总计")
$oldQ3Sheet = $wb.Worksheets.Item("2022-Q3")

$q4Sheet = $wb.Worksheets.Add($oldQ3Sheet)
$q4Sheet.Name = "2022-Q4"

# Reuse the header / index-column formatting already used on the other
# per-quarter sheets (style index carrying bold+border+center alignment)
# by copying it over via PasteSpecial(formats-only), then overwrite the
# values on top so the style is preserved. Re-fetch the source sheet by
# name (rather than reusing the handle captured before the Add() call)
# since sheet handles can go stale once the tab collection shifts.
$q3ForHeader = $wb.Worksheets.Item("2022-Q3")
$q3ForHeader.Range("B1:H1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)

$q3ForIndex = $wb.Worksheets.Item("2022-Q3")
$q3ForIndex.Range("A2:A5").Copy()
$q4Sheet.Range("A2:A5").PasteSpecial(-4122)

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4Sheet.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# row data: index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$data = @(
    @(0, "501062", "南方瑞合三年定期开放混合（LOF）", "7.50", "89.25", "4.17", "0.3128", 6),
    @(1, "001154", "北信瑞丰平安中国主题灵活配置混合", "0.13", "93.42", "3.17", "0.0041", 9),
    @(2, "014668", "银华专精特新量化优选股票A",        "0.26", "94.19", "1.32", "0.0034", 9),
    @(3, "014669", "银华专精特新量化优选股票C",        "0.15", "94.19", "1.32", "0.0020", 9)
)

$row = 2
foreach ($rec in $data) {
    $q4Sheet.Cells.Item($row, 1).Value = $rec[0]

    # Columns B-G hold text values (fund codes / figures stored as strings,
    # matching the other quarter sheets), even though several look numeric.
    # Forcing NumberFormat="@" keeps them as text; ClearFormats() afterwards
    # drops the now-redundant explicit style so the cell matches the plain
    # (unstyled) text cells used elsewhere in the workbook.
    for ($col = 2; $col -le 7; $col++) {
        $cell = $q4Sheet.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $rec[$col - 1]
        $cell.ClearFormats()
    }

    $q4Sheet.Cells.Item($row, 8).Value = $rec[7]
    $row++
}

# ---------------------------------------------------------------------------
# 2) Insert a new row 2 into "总计" for the 2022-Q4 summary figures, shifting
#    the existing rows down by one and renumbering the index column (A).
# ---------------------------------------------------------------------------
$summarySheet.Rows.Item(2).Insert()

# Row 2 keeps the plain (unstyled) look of the other data rows: the Insert()
# above copies the header row's style down onto the new row, so clear that
# before writing the real values.
$summarySheet.Range("B2:D2").ClearFormats()

$summarySheet.Range("B2").Value = "2022-Q4"
$summarySheet.Range("C2").Value = 4
$summarySheet.Range("D2").Value = 0.32

# Renumber the A column (0-based running index) for every data row.
$lastRow = $summarySheet.Cells.Item(1, 1).CurrentRegion.Rows.Count
if ($lastRow -lt 9) { $lastRow = 9 }
for ($r = 2; $r -le $lastRow; $r++) {
    $summarySheet.Cells.Item($r, 1).Value = $r - 2
}

# The freshly-inserted A2 cell has no style yet; give it the same style as
# the rest of the index column by copying formats from A3.
$summarySheetAgain = $wb.Worksheets.Item("总计")
$summarySheetAgain.Range("A3").Copy()
$summarySheetAgain.Range("A2").PasteSpecial(-4122)
